# Auto-generated update of cached market-price/profit columns (H:N) in the
# per-job Leve profit tables, based on refreshed Universalis price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1929.2413
$ws.Range("I40").Value = 1775.0454
$ws.Range("J40").Value = 2413.8572
$ws.Range("K40").Value = 1775.0454
$ws.Range("L40").Value = 2413.8572
$ws.Range("M40").Value = -1600.0454
$ws.Range("N40").Value = -2763.8572

$ws.Range("H98").Value = 801.25
$ws.Range("I98").Value = 801.25
$ws.Range("K98").Value = 801.25
$ws.Range("M98").Value = 696.75

$ws.Range("H122").Value = 801.25
$ws.Range("I122").Value = 801.25
$ws.Range("K122").Value = 2403.75
$ws.Range("M122").Value = 46.25

$ws.Range("H129").Value = 876.0714
$ws.Range("J129").Value = 887.375
$ws.Range("L129").Value = 2662.125
$ws.Range("N129").Value = -12662.125

$ws.Range("H132").Value = 7149606.5
$ws.Range("I132").Value = 8936205
$ws.Range("K132").Value = 26808615
$ws.Range("M132").Value = -26806085

$ws.Range("H137").Value = 1422.1818
$ws.Range("I137").Value = 1428.2307
$ws.Range("K137").Value = 4284.6921
$ws.Range("M137").Value = -1734.6921

$ws.Range("H138").Value = 2152.2415
$ws.Range("I138").Value = 1435.0741
$ws.Range("J138").Value = 2474.9666
$ws.Range("K138").Value = 4305.2223
$ws.Range("L138").Value = 7424.899800000001
$ws.Range("M138").Value = 834.7776999999996
$ws.Range("N138").Value = -17704.8998


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19545.605
$ws.Range("I32").Value = 3520.0833
$ws.Range("K32").Value = 3520.0833
$ws.Range("M32").Value = -3233.0833

$ws.Range("H61").Value = 1035.4
$ws.Range("I61").Value = 849.3461
$ws.Range("J61").Value = 2244.75
$ws.Range("K61").Value = 849.3461
$ws.Range("L61").Value = 2244.75
$ws.Range("M61").Value = -637.3461
$ws.Range("N61").Value = -2668.75

$ws.Range("H74").Value = 3431.85
$ws.Range("I74").Value = 2095.8572
$ws.Range("J74").Value = 6549.1665
$ws.Range("K74").Value = 2095.8572
$ws.Range("L74").Value = 6549.1665
$ws.Range("M74").Value = -1221.8572
$ws.Range("N74").Value = -8297.166499999999

$ws.Range("H77").Value = 3431.85
$ws.Range("I77").Value = 2095.8572
$ws.Range("J77").Value = 6549.1665
$ws.Range("K77").Value = 10479.286
$ws.Range("L77").Value = 32745.8325
$ws.Range("M77").Value = -6111.286
$ws.Range("N77").Value = -41481.8325

$ws.Range("H136").Value = 1035.4
$ws.Range("I136").Value = 849.3461
$ws.Range("J136").Value = 2244.75
$ws.Range("K136").Value = 2548.0383
$ws.Range("L136").Value = 6734.25
$ws.Range("M136").Value = 1.961699999999837
$ws.Range("N136").Value = -11834.25


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2386.5356
$ws.Range("I134").Value = 2386.5356
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7159.6068
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4624.6068
$ws.Range("N134").Value = $null

$ws.Range("H139").Value = 9250
$ws.Range("I139").Value = 10000
$ws.Range("J139").Value = 9000
$ws.Range("K139").Value = 10000
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -4860
$ws.Range("N139").Value = -19280


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 11762.474
$ws.Range("I5").Value = 1168.2858
$ws.Range("J5").Value = 17942.416
$ws.Range("K5").Value = 3504.8574
$ws.Range("L5").Value = 53827.24800000001
$ws.Range("M5").Value = -3392.8574
$ws.Range("N5").Value = -54051.24800000001

$ws.Range("H12").Value = 31.8
$ws.Range("J12").Value = 34.72222
$ws.Range("L12").Value = 104.16666
$ws.Range("N12").Value = -450.16666

$ws.Range("H80").Value = 18751
$ws.Range("J80").Value = 18751
$ws.Range("L80").Value = 56253
$ws.Range("N80").Value = -58125

$ws.Range("H83").Value = 18751
$ws.Range("J83").Value = 18751
$ws.Range("L83").Value = 168759
$ws.Range("N83").Value = -178119

$ws.Range("H122").Value = 3999.9644
$ws.Range("I122").Value = 338.9091
$ws.Range("J122").Value = 17423.834
$ws.Range("K122").Value = 3050.1819
$ws.Range("L122").Value = 156814.506
$ws.Range("M122").Value = -600.1819
$ws.Range("N122").Value = -161714.506

$ws.Range("H131").Value = 801.98
$ws.Range("I131").Value = 407.26666
$ws.Range("J131").Value = 871.6353
$ws.Range("K131").Value = 1221.79998
$ws.Range("L131").Value = 2614.9059
$ws.Range("M131").Value = 3818.20002
$ws.Range("N131").Value = -12694.9059

$ws.Range("H132").Value = 3249.4375
$ws.Range("I132").Value = 2614.6924
$ws.Range("K132").Value = 23532.2316
$ws.Range("M132").Value = -21002.2316

$ws.Range("H135").Value = 11762.474
$ws.Range("I135").Value = 1168.2858
$ws.Range("J135").Value = 17942.416
$ws.Range("K135").Value = 10514.5722
$ws.Range("L135").Value = 161481.744
$ws.Range("M135").Value = -7979.572200000001
$ws.Range("N135").Value = -166551.744


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 13714.667
$ws.Range("J48").Value = 13714.667
$ws.Range("L48").Value = 13714.667
$ws.Range("N48").Value = -14684.667

$ws.Range("H126").Value = 2676429
$ws.Range("I126").Value = 3236.5
$ws.Range("J126").Value = 4203967.5
$ws.Range("K126").Value = 9709.5
$ws.Range("L126").Value = 12611902.5
$ws.Range("M126").Value = -7239.5
$ws.Range("N126").Value = -12616842.5

$ws.Range("H132").Value = 2399.5652
$ws.Range("I132").Value = 1722.2354
$ws.Range("K132").Value = 5166.706200000001
$ws.Range("M132").Value = -2636.706200000001


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = $null

$ws.Range("H22").Value = 737.5454999999999
$ws.Range("I22").Value = 671.3333
$ws.Range("J22").Value = 762.375
$ws.Range("K22").Value = 671.3333
$ws.Range("L22").Value = 762.375
$ws.Range("M22").Value = -376.3333
$ws.Range("N22").Value = -1352.375

$ws.Range("H27").Value = 737.5454999999999
$ws.Range("I27").Value = 671.3333
$ws.Range("J27").Value = 762.375
$ws.Range("K27").Value = 671.3333
$ws.Range("L27").Value = 762.375
$ws.Range("M27").Value = -564.3333
$ws.Range("N27").Value = -976.375

$ws.Range("H100").Value = 1843.25
$ws.Range("J100").Value = 1857.6666
$ws.Range("L100").Value = 1857.6666
$ws.Range("N100").Value = -2939.6666

$ws.Range("H136").Value = 1933.3334
$ws.Range("I136").Value = 1928.5714
$ws.Range("K136").Value = 5785.7142
$ws.Range("M136").Value = -3235.7142


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 250747.5
$ws.Range("I100").Value = 333996.66
$ws.Range("K100").Value = 667993.3199999999
$ws.Range("M100").Value = -667452.3199999999
